$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.485259333333333
$ws.Range("N2").Value = 4.455778
$ws.Range("O2").Value = 0.3057455162066235
$ws.Range("P2").Value = 0.3057455162066235
$ws.Range("Q2").Value = 6.35414983973889
$ws.Range("R2").Value = 57.18734855765
$ws.Range("S2").Value = 0.295448774468984
$ws.Range("T2").Value = 0.295448774468984

# Row 3
$ws.Range("O3").Value = 0.2805555239151429
$ws.Range("P3").Value = 0.2805555239151429
$ws.Range("S3").Value = 0.2711071179052569
$ws.Range("T3").Value = 0.2711071179052569

# Row 4
$ws.Range("O4").Value = 0.4136989598782336
$ws.Range("P4").Value = 0.4136989598782336
$ws.Range("S4").Value = 0.3997666170597783
$ws.Range("T4").Value = 0.3997666170597783

# Row 5
$ws.Range("M5").Value = 1.485259333333333
$ws.Range("N5").Value = 4.455778
$ws.Range("O5").Value = 0.3057455162066235
$ws.Range("P5").Value = 0.3057455162066235
$ws.Range("Q5").Value = 0.2214496911677778
$ws.Range("R5").Value = 1.99304722051
$ws.Range("S5").Value = 0.01029674173763953
$ws.Range("T5").Value = 0.01029674173763953

# Row 6
$ws.Range("O6").Value = 0.2805555239151429
$ws.Range("P6").Value = 0.2805555239151429
$ws.Range("S6").Value = 0.00944840600988606
$ws.Range("T6").Value = 0.009448406009886058

# Row 7
$ws.Range("O7").Value = 0.4136989598782336
$ws.Range("P7").Value = 0.4136989598782336
$ws.Range("S7").Value = 0.01393234281845533
$ws.Range("T7").Value = 0.01393234281845533
